$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 2,3) {
    $ws.Cells.Item($r, 4).Value = -0.128                  # D
    $ws.Cells.Item($r, 11).Value = -24.6                  # K
    $ws.Cells.Item($r, 12).Value = -0.2925089179548157    # L
    $ws.Cells.Item($r, 13).Value = 0.918                  # M
    $ws.Cells.Item($r, 14).Value = 0.02004366812227074    # N
    $ws.Cells.Item($r, 15).Value = -0.03731707317073171   # O
    $ws.Cells.Item($r, 19).Value = 0.918                  # S
    $ws.Cells.Item($r, 20).Value = 1                      # T
    $ws.Cells.Item($r, 21).Value = 157.3                  # U
    $ws.Cells.Item($r, 22).Value = 3.434497816593887      # V
    $ws.Cells.Item($r, 23).Value = -0.05023483765570758   # W
    $ws.Cells.Item($r, 24).Value = 0.3414318517837213     # X
    $ws.Cells.Item($r, 25).Value = -0.3916666894394288    # Y
    $ws.Cells.Item($r, 26).Value = 0.1074623051367237     # Z
    $ws.Cells.Item($r, 28).Value = 0.07543545080714474    # AB
    $ws.Cells.Item($r, 29).Value = -0.07543545080714474   # AC
    $ws.Cells.Item($r, 30).Value = 420.9                  # AD
    $ws.Cells.Item($r, 32).Value = 420.9                  # AF
    $ws.Cells.Item($r, 33).Value = 263.6                  # AG
    $ws.Cells.Item($r, 34).Value = 0.9018641525605313     # AH
    $ws.Cells.Item($r, 35).Value = 0.5178395669291339     # AI
    $ws.Cells.Item($r, 36).Value = 0.8519715578539108     # AJ
    $ws.Cells.Item($r, 37).Value = 0.4021357742181541     # AK
}
